$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# add_new_qualification_types: just a cursor/selection move
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("add_new_qualification_types")
$ws.Range("C19").Select()

# ---------------------------------------------------------------------
# edit_qualification_types: just a cursor/selection move
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("edit_qualification_types")
$ws.Range("C2").Select()

# ---------------------------------------------------------------------
# add_new_entry_criteria: trimmed trailing formatting-only rows/cols,
# moved the selection to column K (selecting the whole column)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("add_new_entry_criteria")
$ws.Rows.Item(8).EntireRow.Delete()
$ws.Rows.Item(8).EntireRow.Delete()
$ws.Rows.Item(8).EntireRow.Delete()
$ws.Columns.Item(11).ColumnWidth = 15.28515625
$ws.Columns.Item(12).ColumnWidth = 15.28515625
$ws.Columns.Item(13).ColumnWidth = 15.28515625
$ws.Columns.Item(11).Select()

# ---------------------------------------------------------------------
# search_entry_criteria: just a cursor/selection move
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("search_entry_criteria")
$ws.Range("J18").Select()

# ---------------------------------------------------------------------
# edit_entry_criteria: new "type" / "outcome" columns added, code/name
# swapped between the two data rows, new type + outcome values filled
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("edit_entry_criteria")
$ws.Range("D1:E1").EntireColumn.Insert()

$ws.Range("D1").Value = "new type"
$ws.Range("E1").Value = "new outcome"
$ws.Range("G1").Value = "new grade"
$ws.Range("H1").Value = "new value"

$ws.Range("B2").Value = "Auto_UPD_EC_001"
$ws.Range("C2").Value = "Auto_Sri Lankan A/L"
$ws.Range("D2").Value = "Sri Lankan A/L"
$ws.Range("E2").Value = "bio"

$ws.Range("B3").Value = "Auto_UPD_EC_002"
$ws.Range("C3").Value = "Auto_UPD London A/L"
$ws.Range("D3").Value = "London A/L"
$ws.Range("E3").Value = "bio"

$ws.Columns.Item(2).ColumnWidth = 17
$ws.Columns.Item(4).ColumnWidth = 13.42578125
$ws.Columns.Item(5).ColumnWidth = 16.42578125
$ws.Columns.Item(6).ColumnWidth = 19

$ws.Activate()
$ws.Range("D4").Select()

# ---------------------------------------------------------------------
# filter_programs: no longer the active tab (edit_entry_criteria is)
# handled implicitly by activating edit_entry_criteria above, which
# Excel records as the new tabSelected sheet.
# ---------------------------------------------------------------------
